# Hardik Pandya - IPL Mumbai Indians innings-by-innings log.
# The sheet originally held only one match row (row 2); this adds 11 more
# match rows, expanding the table to A1:K13 (header + 12 innings).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: push the existing row 2 down to row 3, and open up 9 more blank
# rows below it (rows 4-12) so every new record (rows 2-13) has a home.
$ws.Rows("2:2").Insert()
$ws.Rows("4:12").Insert()

# Row 2: Hardik Pandya vs Kolkata Knight Riders (Oct 16 2020)
$ws.Range("G2:K2").NumberFormat = "@"
$ws.Cells.Item(2,1).Value = " Oct 16 2020"
$ws.Cells.Item(2,2).Value = " Abu Dhabi"
$ws.Cells.Item(2,3).Value = "Mumbai won by 8 wickets (with 19 balls remaining)"
$ws.Cells.Item(2,4).Value = "Mumbai Indians"
$ws.Cells.Item(2,5).Value = "Kolkata Knight Riders"
$ws.Cells.Item(2,6).Value = "Hardik Pandya "
$ws.Cells.Item(2,7).Value = "21"
$ws.Cells.Item(2,8).Value = "11"
$ws.Cells.Item(2,9).Value = "3"
$ws.Cells.Item(2,10).Value = "1"
$ws.Cells.Item(2,11).Value = "190.90"
$ws.Range("G2:K2").Style = "Normal"

# Row 3: Hardik Pandya vs Delhi Capitals (Nov 10 2020)
$ws.Range("G3:K3").NumberFormat = "@"
$ws.Cells.Item(3,1).Value = " Nov 10 2020"
$ws.Cells.Item(3,2).Value = " Dubai (DSC)"
$ws.Cells.Item(3,3).Value = "Mumbai won by 5 wickets (with 8 balls remaining)"
$ws.Cells.Item(3,4).Value = "Mumbai Indians"
$ws.Cells.Item(3,5).Value = "Delhi Capitals"
$ws.Cells.Item(3,6).Value = "Hardik Pandya "
$ws.Cells.Item(3,7).Value = "3"
$ws.Cells.Item(3,8).Value = "5"
$ws.Cells.Item(3,9).Value = "0"
$ws.Cells.Item(3,10).Value = "0"
$ws.Cells.Item(3,11).Value = "60.00"
$ws.Range("G3:K3").Style = "Normal"

# Row 4: Hardik Pandya vs Rajasthan Royals (Oct 6 2020)
$ws.Range("G4:K4").NumberFormat = "@"
$ws.Cells.Item(4,1).Value = " Oct 6 2020"
$ws.Cells.Item(4,2).Value = " Abu Dhabi"
$ws.Cells.Item(4,3).Value = "Mumbai won by 57 runs"
$ws.Cells.Item(4,4).Value = "Mumbai Indians"
$ws.Cells.Item(4,5).Value = "Rajasthan Royals"
$ws.Cells.Item(4,6).Value = "Hardik Pandya "
$ws.Cells.Item(4,7).Value = "30"
$ws.Cells.Item(4,8).Value = "19"
$ws.Cells.Item(4,9).Value = "2"
$ws.Cells.Item(4,10).Value = "1"
$ws.Cells.Item(4,11).Value = "157.89"
$ws.Range("G4:K4").Style = "Normal"

# Row 5: Hardik Pandya vs Kings XI Punjab (Oct 18 2020)
$ws.Range("G5:K5").NumberFormat = "@"
$ws.Cells.Item(5,1).Value = " Oct 18 2020"
$ws.Cells.Item(5,2).Value = " Dubai (DSC)"
$ws.Cells.Item(5,3).Value = "Match tied (Kings XI won the one-over eliminator)"
$ws.Cells.Item(5,4).Value = "Mumbai Indians"
$ws.Cells.Item(5,5).Value = "Kings XI Punjab"
$ws.Cells.Item(5,6).Value = "Hardik Pandya "
$ws.Cells.Item(5,7).Value = "8"
$ws.Cells.Item(5,8).Value = "4"
$ws.Cells.Item(5,9).Value = "0"
$ws.Cells.Item(5,10).Value = "1"
$ws.Cells.Item(5,11).Value = "200.00"
$ws.Range("G5:K5").Style = "Normal"

# Row 6: Hardik Pandya vs Rajasthan Royals (Oct 25 2020)
$ws.Range("G6:K6").NumberFormat = "@"
$ws.Cells.Item(6,1).Value = " Oct 25 2020"
$ws.Cells.Item(6,2).Value = " Abu Dhabi"
$ws.Cells.Item(6,3).Value = "Royals won by 8 wickets (with 10 balls remaining)"
$ws.Cells.Item(6,4).Value = "Mumbai Indians"
$ws.Cells.Item(6,5).Value = "Rajasthan Royals"
$ws.Cells.Item(6,6).Value = "Hardik Pandya "
$ws.Cells.Item(6,7).Value = "60"
$ws.Cells.Item(6,8).Value = "21"
$ws.Cells.Item(6,9).Value = "2"
$ws.Cells.Item(6,10).Value = "7"
$ws.Cells.Item(6,11).Value = "285.71"
$ws.Range("G6:K6").Style = "Normal"

# Row 7: Hardik Pandya vs Sunrisers Hyderabad (Oct 4 2020)
$ws.Range("G7:K7").NumberFormat = "@"
$ws.Cells.Item(7,1).Value = " Oct 4 2020"
$ws.Cells.Item(7,2).Value = " Sharjah"
$ws.Cells.Item(7,3).Value = "Mumbai won by 34 runs"
$ws.Cells.Item(7,4).Value = "Mumbai Indians"
$ws.Cells.Item(7,5).Value = "Sunrisers Hyderabad"
$ws.Cells.Item(7,6).Value = "Hardik Pandya "
$ws.Cells.Item(7,7).Value = "28"
$ws.Cells.Item(7,8).Value = "19"
$ws.Cells.Item(7,9).Value = "2"
$ws.Cells.Item(7,10).Value = "2"
$ws.Cells.Item(7,11).Value = "147.36"
$ws.Range("G7:K7").Style = "Normal"

# Row 8: Hardik Pandya vs Delhi Capitals (Oct 11 2020)
$ws.Range("G8:K8").NumberFormat = "@"
$ws.Cells.Item(8,1).Value = " Oct 11 2020"
$ws.Cells.Item(8,2).Value = " Abu Dhabi"
$ws.Cells.Item(8,3).Value = "Mumbai won by 5 wickets (with 2 balls remaining)"
$ws.Cells.Item(8,4).Value = "Mumbai Indians"
$ws.Cells.Item(8,5).Value = "Delhi Capitals"
$ws.Cells.Item(8,6).Value = "Hardik Pandya "
$ws.Cells.Item(8,7).Value = "0"
$ws.Cells.Item(8,8).Value = "2"
$ws.Cells.Item(8,9).Value = "0"
$ws.Cells.Item(8,10).Value = "0"
$ws.Cells.Item(8,11).Value = "0.00"
$ws.Range("G8:K8").Style = "Normal"

# Row 9: Hardik Pandya vs Delhi Capitals (Nov 5 2020)
$ws.Range("G9:K9").NumberFormat = "@"
$ws.Cells.Item(9,1).Value = " Nov 5 2020"
$ws.Cells.Item(9,2).Value = " Dubai (DSC)"
$ws.Cells.Item(9,3).Value = "Mumbai won by 57 runs"
$ws.Cells.Item(9,4).Value = "Mumbai Indians"
$ws.Cells.Item(9,5).Value = "Delhi Capitals"
$ws.Cells.Item(9,6).Value = "Hardik Pandya "
$ws.Cells.Item(9,7).Value = "37"
$ws.Cells.Item(9,8).Value = "14"
$ws.Cells.Item(9,9).Value = "0"
$ws.Cells.Item(9,10).Value = "5"
$ws.Cells.Item(9,11).Value = "264.28"
$ws.Range("G9:K9").Style = "Normal"

# Row 10: Hardik Pandya vs Chennai Super Kings (Sep 19 2020)
$ws.Range("G10:K10").NumberFormat = "@"
$ws.Cells.Item(10,1).Value = " Sep 19 2020"
$ws.Cells.Item(10,2).Value = " Abu Dhabi"
$ws.Cells.Item(10,3).Value = "Super Kings won by 5 wickets (with 4 balls remaining)"
$ws.Cells.Item(10,4).Value = "Mumbai Indians"
$ws.Cells.Item(10,5).Value = "Chennai Super Kings"
$ws.Cells.Item(10,6).Value = "Hardik Pandya "
$ws.Cells.Item(10,7).Value = "14"
$ws.Cells.Item(10,8).Value = "10"
$ws.Cells.Item(10,9).Value = "0"
$ws.Cells.Item(10,10).Value = "2"
$ws.Cells.Item(10,11).Value = "140.00"
$ws.Range("G10:K10").Style = "Normal"

# Row 11: Hardik Pandya vs Royal Challengers Bangalore (Sep 28 2020)
$ws.Range("G11:K11").NumberFormat = "@"
$ws.Cells.Item(11,1).Value = " Sep 28 2020"
$ws.Cells.Item(11,2).Value = " Dubai (DSC)"
$ws.Cells.Item(11,3).Value = "Match tied (RCB won the one-over eliminator)"
$ws.Cells.Item(11,4).Value = "Mumbai Indians"
$ws.Cells.Item(11,5).Value = "Royal Challengers Bangalore"
$ws.Cells.Item(11,6).Value = "Hardik Pandya "
$ws.Cells.Item(11,7).Value = "15"
$ws.Cells.Item(11,8).Value = "13"
$ws.Cells.Item(11,9).Value = "0"
$ws.Cells.Item(11,10).Value = "1"
$ws.Cells.Item(11,11).Value = "115.38"
$ws.Range("G11:K11").Style = "Normal"

# Row 12: Hardik Pandya vs Kings XI Punjab (Oct 1 2020)
$ws.Range("G12:K12").NumberFormat = "@"
$ws.Cells.Item(12,1).Value = " Oct 1 2020"
$ws.Cells.Item(12,2).Value = " Abu Dhabi"
$ws.Cells.Item(12,3).Value = "Mumbai won by 48 runs"
$ws.Cells.Item(12,4).Value = "Mumbai Indians"
$ws.Cells.Item(12,5).Value = "Kings XI Punjab"
$ws.Cells.Item(12,6).Value = "Hardik Pandya "
$ws.Cells.Item(12,7).Value = "30"
$ws.Cells.Item(12,8).Value = "11"
$ws.Cells.Item(12,9).Value = "3"
$ws.Cells.Item(12,10).Value = "2"
$ws.Cells.Item(12,11).Value = "272.72"
$ws.Range("G12:K12").Style = "Normal"

# Row 13: Hardik Pandya vs Kolkata Knight Riders (Sep 23 2020)
$ws.Range("G13:K13").NumberFormat = "@"
$ws.Cells.Item(13,1).Value = " Sep 23 2020"
$ws.Cells.Item(13,2).Value = " Abu Dhabi"
$ws.Cells.Item(13,3).Value = "Mumbai won by 49 runs"
$ws.Cells.Item(13,4).Value = "Mumbai Indians"
$ws.Cells.Item(13,5).Value = "Kolkata Knight Riders"
$ws.Cells.Item(13,6).Value = "Hardik Pandya "
$ws.Cells.Item(13,7).Value = "18"
$ws.Cells.Item(13,8).Value = "13"
$ws.Cells.Item(13,9).Value = "2"
$ws.Cells.Item(13,10).Value = "1"
$ws.Cells.Item(13,11).Value = "138.46"
$ws.Range("G13:K13").Style = "Normal"
